# Leave Card update for ANCIRO, DANILO workbook
# - Switches the Column A "effectivity date" ladder on Sheet1 (the big leave
#   table) from EDATE(...) to EOMONTH(...), row by row (each cell gets its
#   own explicit formula instead of a shared formula group).
# - Corrects the final manual entry date (A293) by one day.
# - Updates the "OPTIONAL RETIREMENT EFFECTIVE DATE" note text.
# - Leaves "Sheet1" (the big table) as the active/selected tab instead of
#   "CONVERTION".

$wb = $excel.ActiveWorkbook

# --- Sheet1 (physically sheet2.xml) : the big leave-card table -------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()

$ws2.Range("A13").Formula = "=EOMONTH(A12,1)"
$ws2.Range("A14:A20").Formula = "=EOMONTH(A13,1)"
$ws2.Range("A22").Formula = "=EOMONTH(A20,1)"
$ws2.Range("A23:A33").Formula = "=EOMONTH(A22,1)"
$ws2.Range("A35").Formula = "=EOMONTH(A33,1)"
$ws2.Range("A36:A46").Formula = "=EOMONTH(A35,1)"
$ws2.Range("A48").Formula = "=EOMONTH(A46,1)"
$ws2.Range("A49:A59").Formula = "=EOMONTH(A48,1)"
$ws2.Range("A61").Formula = "=EOMONTH(A59,1)"
$ws2.Range("A62:A72").Formula = "=EOMONTH(A61,1)"
$ws2.Range("A74").Formula = "=EOMONTH(A72,1)"
$ws2.Range("A75:A85").Formula = "=EOMONTH(A74,1)"
$ws2.Range("A87").Formula = "=EOMONTH(A85,1)"
$ws2.Range("A88:A98").Formula = "=EOMONTH(A87,1)"
$ws2.Range("A100").Formula = "=EOMONTH(A98,1)"
$ws2.Range("A101:A111").Formula = "=EOMONTH(A100,1)"
$ws2.Range("A113").Formula = "=EOMONTH(A111,1)"
$ws2.Range("A114:A124").Formula = "=EOMONTH(A113,1)"
$ws2.Range("A126").Formula = "=EOMONTH(A124,1)"
$ws2.Range("A127:A137").Formula = "=EOMONTH(A126,1)"
$ws2.Range("A139").Formula = "=EOMONTH(A137,1)"
$ws2.Range("A140:A150").Formula = "=EOMONTH(A139,1)"
$ws2.Range("A152").Formula = "=EOMONTH(A150,1)"
$ws2.Range("A153:A163").Formula = "=EOMONTH(A152,1)"
$ws2.Range("A165").Formula = "=EOMONTH(A163,1)"
$ws2.Range("A166:A176").Formula = "=EOMONTH(A165,1)"
$ws2.Range("A178").Formula = "=EOMONTH(A176,1)"
$ws2.Range("A179:A189").Formula = "=EOMONTH(A178,1)"
$ws2.Range("A191").Formula = "=EOMONTH(A189,1)"
$ws2.Range("A192:A202").Formula = "=EOMONTH(A191,1)"
$ws2.Range("A204").Formula = "=EOMONTH(A202,1)"
$ws2.Range("A205:A215").Formula = "=EOMONTH(A204,1)"
$ws2.Range("A217").Formula = "=EOMONTH(A215,1)"
$ws2.Range("A218:A228").Formula = "=EOMONTH(A217,1)"
$ws2.Range("A230").Formula = "=EOMONTH(A228,1)"
$ws2.Range("A231:A241").Formula = "=EOMONTH(A230,1)"
$ws2.Range("A243").Formula = "=EOMONTH(A241,1)"
$ws2.Range("A244:A254").Formula = "=EOMONTH(A243,1)"
$ws2.Range("A256").Formula = "=EOMONTH(A254,1)"
$ws2.Range("A257:A267").Formula = "=EOMONTH(A256,1)"
$ws2.Range("A269").Formula = "=EOMONTH(A267,1)"
$ws2.Range("A270:A280").Formula = "=EOMONTH(A269,1)"
$ws2.Range("A282").Formula = "=EOMONTH(A280,1)"
$ws2.Range("A283:A292").Formula = "=EOMONTH(A282,1)"
# Final manually-entered retirement date, one day earlier.
$ws2.Range("A293").Value = 44907

# Retirement note text (day bumped from the 12th to the 13th).
$ws2.Range("B294").Value = "OPTIONAL RETIREMENT EFFECTIVE DATE: DECEMBER 13, 2022"

# Scroll/selection bookkeeping on Sheet1: default-pane selection on the
# header block, and the frozen/split bottom pane parked near the new last
# row of data.
$ws2.Range("B2:C2").Select()
$ws2.Range("H296").Select()

# --- CONVERTION sheet (sheet3.xml) : no longer the active tab --------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Range("J20").Select()

# --- Make Sheet1 the active/selected tab again ------------------------------
$ws2.Activate()
